$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 updates ---
$ws.Range("D2").Value = "Design UI"
$ws.Range("E2").Value = 100
$ws.Range("F2").ClearContents()
$ws.Range("H2").Value = 50
$ws.Range("I2").Value = 1234

# --- Row 3 updates ---
$ws.Range("A3").Value = 2
$ws.Range("D3").Value = "Implement Backend"
$ws.Range("E3").Value = 55
$ws.Range("H3").ClearContents()
$ws.Range("I3").Value = 1235

# --- Row 4 updates ---
$ws.Range("D4").Value = "Testing"
$ws.Range("E4").Value = 15
$ws.Range("I4").Value = 123

# --- Row 5 updates ---
$ws.Range("F5").ClearContents()

# --- Row 9 updates ---
$ws.Range("D9").Value = "Analysis"
$ws.Range("I9").Value = 3135

# --- New rows 10-13 ---
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = 4
$ws.Range("C10").Value = "Project D"
$ws.Range("D10").Value = "Database Setup"
$ws.Range("E10").Value = 75
$ws.Range("I10").Value = 1314

$ws.Range("A11").Value = 10
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = "Project D"
$ws.Range("D11").Value = "Implement Backend"
$ws.Range("E11").Value = 47
$ws.Range("I11").Value = 13486

$ws.Range("A12").Value = 11
$ws.Range("B12").Value = 4
$ws.Range("C12").Value = "Project D"
$ws.Range("D12").Value = "UI Development"
$ws.Range("E12").Value = 46
$ws.Range("I12").Value = 1314

$ws.Range("A13").Value = 12
$ws.Range("B13").Value = 4
$ws.Range("C13").Value = "Project D"
$ws.Range("D13").Value = "Testing"
$ws.Range("E13").Value = 47
$ws.Range("I13").Value = 13486

# --- Resize the table to cover the new rows ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:I13"))

# --- Update the sheet view (zoom / scroll / selection) ---
$excel.ActiveWindow.Zoom = 85
$ws.Range("A1:I13").Select()
